$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1. Insert a new row above row 167 by duplicating it (keeps formats/styles identical
#    to the existing row167, matching how the real edit preserved A167's style/value).
$ws.Rows.Item(167).Copy()
$ws.Rows.Item(167).Insert()

# 2. Update the newly inserted row167 with the corrected/new log entry.
$ws.Range("B167").Value2 = "Tag"
$ws.Range("D167").Value2 = "Navleen Purewal"
$ws.Range("D167").Style = $ws.Range("C167").Style
$ws.Range("E167").Value2 = 45793
$ws.Range("F167").Value2 = 0.39097222222222222
$ws.Range("G167").Formula = "=((E167+F167)-DATE(1970,1,1))*86400"

# 3. Append two brand-new rows at the bottom of the table for the missed log entries.
$newRow1 = $lo.ListRows.Add()
$newRow1.Range.Cells(1,1).Value2 = 173
$newRow1.Range.Cells(1,2).Value2 = "Tag"
$newRow1.Range.Cells(1,3).Value2 = "Cora Walshe"
$newRow1.Range.Cells(1,4).Value2 = "Keira Kelly"
$newRow1.Range.Cells(1,5).Value2 = 45778
$newRow1.Range.Cells(1,6).Value2 = 0.59444444444444444
$newRow1.Range.Cells(1,7).Formula = "=((E172+F172)-DATE(1970,1,1))*86400"

$newRow2 = $lo.ListRows.Add()
$newRow2.Range.Cells(1,1).Value2 = 174
$newRow2.Range.Cells(1,2).Value2 = "Tag"
$newRow2.Range.Cells(1,3).Value2 = "Samantha Crowder"
$newRow2.Range.Cells(1,4).Value2 = "Harrison Bartley"
$newRow2.Range.Cells(1,5).Value2 = 45792
$newRow2.Range.Cells(1,6).Value2 = 0.5
$newRow2.Range.Cells(1,7).Formula = "=((E173+F173)-DATE(1970,1,1))*86400"

# 4. Remove the stray "Column1" table column (H) that was never used for data.
$lo.ListColumns.Item(8).Delete()

# 5. Restore the view to where the user left off editing.
$ws.Application.ActiveWindow.ScrollRow = 160
$ws.Range("D174").Select()
